# Update the "employment APS data" period text and move the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B9").Value = "Jul 2022 - Jun 2023 data"

$ws.Range("B2:B9").Select()
